# Insert a new weekly price-report row for "Poroto granado" at the top of the
# existing data block (before the current row 588), pushing all rows from
# 588 downward by one. This mirrors a new consolidated daily/weekly data
# point being added to the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 588; everything at/after 588 shifts down.
$ws.Rows.Item(588).Insert()

# Populate the newly inserted row 588 with the new record's data.
$ws.Cells.Item(588, 1).Value  = 6
$ws.Cells.Item(588, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(588, 3).Value  = "Metropolitana"
$ws.Cells.Item(588, 4).Value  = 44995
$ws.Cells.Item(588, 5).Value  = 13
$ws.Cells.Item(588, 6).Value  = 100112030
$ws.Cells.Item(588, 7).Value  = "Poroto granado"
$ws.Cells.Item(588, 8).Value  = "Sin especificar"
$ws.Cells.Item(588, 9).Value  = "Primera"
$ws.Cells.Item(588, 10).Value = 680
$ws.Cells.Item(588, 11).Value = 30000
$ws.Cells.Item(588, 12).Value = 32000
$ws.Cells.Item(588, 13).Value = 31235
$ws.Cells.Item(588, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(588, 15).Value = "Región Metropolitana"
$ws.Cells.Item(588, 16).Value = 1249
$ws.Cells.Item(588, 17).Value = 25
$ws.Cells.Item(588, 18).Value = "Hortaliza"
